# Fix for having two or more links in normative text
# Adds two new rows ("hyperlink5" and "hyperlink6") to the "Normative Rules"
# table, directly after the existing "hyperlink4" row (and before "table1"),
# each demonstrating multiple norm:superscript / norm:subscript style inline
# anchors inside a single piece of normative text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet has a single Excel Table ("Table1") that currently spans
# A1:F46 (header row + 45 data rows). The new rows need to land right after
# the existing "hyperlink4" row, which is row 24, i.e. they become the new
# rows 25 and 26, pushing every following row down by two.
$lo = $ws.ListObjects.Item(1)

# Insert two blank worksheet rows right before the current row 25 ("table1"),
# shifting it (and everything below it) down by two rows.
$ws.Rows("25:26").Insert()

# Row 25: hyperlink5
$ws.Range("A25").Value = "my-chapter_name"
$ws.Range("B25").Value = "hyperlink5"
$ws.Range("C25").Value = "GHI &lt;&lt;norm:superscript&gt;&gt; and &lt;&lt;norm:subscript&gt;&gt; JKL"
$ws.Range("D25").Value = '["norm:hyperlink5"]'

# Row 26: hyperlink6
$ws.Range("A26").Value = "my-chapter_name"
$ws.Range("B26").Value = "hyperlink6"
$ws.Range("C26").Value = "JKL &lt;&lt;norm:superscript,hello&gt;&gt; and &lt;&lt;norm:subscript,goodbye&gt;&gt; MNO"
$ws.Range("D26").Value = '["norm:hyperlink6"]'

# Grow the table so it once again covers the whole data block, now
# A1:F48 (header + 47 data rows).
$lastRow = $lo.Range.Row + $lo.Range.Rows.Count - 1 + 2
$lo.Resize($ws.Range("A1:F48"))
